$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '89.591.64'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.137.57'
$ws.Range("E3").Value = '  -4.42%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.43'
$ws.Range("E5").Value = '  -1.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '635.59'
$ws.Range("E6").Value = '  +0.71%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.397'
$ws.Range("E7").Value = '  -1.63%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.762'
$ws.Range("E8").Value = '  +6.38%  '

$ws.Range("E9").Value = '  +0.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.135.20'
$ws.Range("E10").Value = '  -4.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.560'
$ws.Range("E11").Value = '  -5.05%  '

$ws.Range("E12").Value = '  -0.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000252'
$ws.Range("E13").Value = '  -5.92%  '

$ws.Range("E14").Value = '  -0.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.237.08'
$ws.Range("E15").Value = '  -3.17%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.715.42'
$ws.Range("E16").Value = '  -4.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '32.48'
$ws.Range("E17").Value = '  -5.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.137.97'
$ws.Range("E18").Value = '  -3.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000228'
$ws.Range("E19").Value = '  +9.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.37'
$ws.Range("E20").Value = '  +1.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.34'
$ws.Range("E21").Value = '  -5.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '428.78'
$ws.Range("E22").Value = '  -3.10%  '

$ws.Range("E23").Value = '  -6.04%  '

$ws.Range("E24").Value = '  -6.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.46'
$ws.Range("E25").Value = '  +1.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '82.32'
$ws.Range("E26").Value = '  +6.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.59'
$ws.Range("E27").Value = '  -5.83%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.295.49'
$ws.Range("E28").Value = '  -4.43%  '

$ws.Range("E29").Value = '  +0.19%  '

$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.996'
$ws.Range("E30").Value = '  -0.15%  '

$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.157'
$ws.Range("E31").Value = '  -10.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.02'
$ws.Range("E32").Value = '  +6.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.24'
$ws.Range("E33").Value = '  -6.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '507.76'
$ws.Range("E34").Value = '  -9.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.148'
$ws.Range("E35").Value = '  +12.85%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.14'
$ws.Range("E36").Value = '  -1.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.30'
$ws.Range("E37").Value = '  +0.65%  '

$ws.Range("E38").Value = '  -4.72%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.06'
$ws.Range("E39").Value = '  -2.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.26'
$ws.Range("E40").Value = '  -0.72%  '

$ws.Range("E41").Value = '  +0.35%  '

$ws.Range("E42").Value = '  -0.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.88'
$ws.Range("E43").Value = '  -6.68%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.365'
$ws.Range("E44").Value = '  -7.60%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '145.99'
$ws.Range("E45").Value = '  -2.96%  '

$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.87'
$ws.Range("E46").Value = '  -3.86%  '

$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.130'
$ws.Range("E47").Value = '  +0.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '164.77'
$ws.Range("E48").Value = '  -8.41%  '

$ws.Range("B49").Value = 'Hedera'
$ws.Range("C49").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0650'
$ws.Range("E49").Value = '  +7.90%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.723'
$ws.Range("E50").Value = '  -0.95%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.36'
$ws.Range("E51").Value = '  -2.87%  '
